$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Meeting end time: "17:30" -> "16:30" (the sprint actually ran until
#    16:30). Word records the last edit position with the reserved
#    "_GoBack" bookmark, which moves here and away from the old edit
#    location later in the document (near "crow's"). Adding the new
#    "_GoBack" bookmark automatically removes the old one.
# ------------------------------------------------------------------
$timeRange = $d.Content
$timeRange.Find.Execute("17:30") | Out-Null
$sevenPos = $timeRange.Start + 1
$d.Range($sevenPos, $sevenPos + 1).Text = "6"
$d.Bookmarks.Add("_GoBack", $d.Range($sevenPos + 1, $sevenPos + 1)) | Out-Null

# ------------------------------------------------------------------
# 2. Merge the split "Create the presentation for the fourth
#    presentation..." sentence back into a single run (re-typed as one
#    continuous piece of text).
# ------------------------------------------------------------------
$firstFourth = $d.Content
$firstFourth.Find.Execute("fourth presentation") | Out-Null
$secondFourth = $d.Range($firstFourth.End, $d.Content.End)
$secondFourth.Find.Execute("fourth presentation") | Out-Null
$presoPara = $secondFourth.Paragraphs(1).Range
$presoRange = $d.Range($presoPara.Start, $presoPara.End - 1)
$presoText = $presoRange.Text
$presoRange.Find.Execute($presoText, $false, $false, $false, $false, $false, $true, 1, $false, $presoText, 2) | Out-Null

# ------------------------------------------------------------------
# 3 & 4. Merge the trailing ", finish implementing task durations (30m)"
#    and ", fix bugs with different task durations (30m)" back into a
#    single run each, without disturbing the separate "studio jam" run
#    before them. A temporary bookmark is used to "fence off" the
#    preceding run while the edit is made, then removed again.
# ------------------------------------------------------------------
function Merge-TrailingClause($fullClause) {
    $found = $d.Content
    $found.Find.Execute($fullClause) | Out-Null
    $commaPos = $found.Start
    $d.Bookmarks.Add("_MergeFence", $d.Range($commaPos, $commaPos)) | Out-Null
    $clauseRange = $d.Range($commaPos, $commaPos + $fullClause.Length)
    $clauseRange.Find.Execute($fullClause, $false, $false, $false, $false, $false, $true, 1, $false, $fullClause, 2) | Out-Null
    $d.Bookmarks("_MergeFence").Delete()
}

Merge-TrailingClause ", finish implementing task durations (30m)"
Merge-TrailingClause ", fix bugs with different task durations (30m)"
